# Update column F (view/attendance counts) values across the 4 worksheets
# to match regenerated stats output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 308
$ws.Range("F6").Value = 979
$ws.Range("F8").Value = 2409
$ws.Range("F10").Value = 1179
$ws.Range("F11").Value = 885
$ws.Range("F12").Value = 588
$ws.Range("F13").Value = 873
$ws.Range("F14").Value = 1065
$ws.Range("F16").Value = 283
$ws.Range("F17").Value = 109
$ws.Range("F18").Value = 704
$ws.Range("F19").Value = 740
$ws.Range("F20").Value = 172
$ws.Range("F21").Value = 449
$ws.Range("F22").Value = 1088
$ws.Range("F23").Value = 73
$ws.Range("F24").Value = 528
$ws.Range("F25").Value = 561
$ws.Range("F26").Value = 204
$ws.Range("F27").Value = 284
$ws.Range("F28").Value = 283
$ws.Range("F29").Value = 658
$ws.Range("F30").Value = 2941
$ws.Range("F31").Value = 452
$ws.Range("F32").Value = 40
$ws.Range("F33").Value = 296
$ws.Range("F36").Value = 112
$ws.Range("F37").Value = 1557
$ws.Range("F39").Value = 140
$ws.Range("F40").Value = 75
$ws.Range("F41").Value = 125
$ws.Range("F42").Value = 59
$ws.Range("F44").Value = 107
$ws.Range("F45").Value = 114
$ws.Range("F46").Value = 75
$ws.Range("F47").Value = 17

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 167
$ws.Range("F13").Value = 18
$ws.Range("F14").Value = 148
$ws.Range("F18").Value = 3

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2246
$ws.Range("F3").Value = 702

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2246
$ws.Range("F3").Value = 702
$ws.Range("F4").Value = 12
$ws.Range("F8").Value = 979
$ws.Range("F9").Value = 2409
$ws.Range("F11").Value = 1179
$ws.Range("F12").Value = 885
$ws.Range("F13").Value = 588
$ws.Range("F14").Value = 873
$ws.Range("F15").Value = 1065
$ws.Range("F16").Value = 283
$ws.Range("F18").Value = 109
$ws.Range("F19").Value = 704
$ws.Range("F21").Value = 740
$ws.Range("F22").Value = 172
$ws.Range("F23").Value = 450
$ws.Range("F24").Value = 1088
$ws.Range("F26").Value = 73
$ws.Range("F27").Value = 528
$ws.Range("F28").Value = 561
$ws.Range("F29").Value = 204
$ws.Range("F30").Value = 283
$ws.Range("F32").Value = 2942
$ws.Range("F33").Value = 167
$ws.Range("F34").Value = 452
$ws.Range("F35").Value = 40
$ws.Range("F37").Value = 112
$ws.Range("F38").Value = 1557
$ws.Range("F41").Value = 18
$ws.Range("F42").Value = 125
$ws.Range("F43").Value = 59
$ws.Range("F45").Value = 107
$ws.Range("F47").Value = 75
$ws.Range("F49").Value = 3

